$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text content could be misread as a numeric value by Excel
# (e.g. "591.55", "1.00") are first forced to Text format so the literal
# string is preserved, then the format is reset back to General so no
# extra style/formatting is left behind on the cell.

$textCells = @(
    @{ Addr = "D5"; Val = "591.55" }
    @{ Addr = "D6"; Val = "147.76" }
    @{ Addr = "D9"; Val = "0.532" }
    @{ Addr = "D11"; Val = "5.91" }
    @{ Addr = "D12"; Val = "0.464" }
    @{ Addr = "D13"; Val = "0.0000249" }
    @{ Addr = "D14"; Val = "37.31" }
    @{ Addr = "D17"; Val = "7.22" }
    @{ Addr = "D20"; Val = "468.93" }
    @{ Addr = "D21"; Val = "14.45" }
    @{ Addr = "D22"; Val = "0.739" }
    @{ Addr = "D23"; Val = "7.48" }
    @{ Addr = "D24"; Val = "2.37" }
    @{ Addr = "D25"; Val = "13.10" }
    @{ Addr = "D26"; Val = "81.26" }
    @{ Addr = "D28"; Val = "9.79" }
    @{ Addr = "D29"; Val = "2.72" }
    @{ Addr = "D30"; Val = "2.23" }
    @{ Addr = "D31"; Val = "7.29" }
    @{ Addr = "D32"; Val = "1.00" }
    @{ Addr = "D33"; Val = "27.60" }
    @{ Addr = "D34"; Val = "0.110" }
    @{ Addr = "D37"; Val = "2.36" }
    @{ Addr = "D38"; Val = "6.14" }
    @{ Addr = "D40"; Val = "51.95" }
    @{ Addr = "D41"; Val = "456.99" }
    @{ Addr = "D42"; Val = "9.16" }
    @{ Addr = "D43"; Val = "0.294" }
    @{ Addr = "D46"; Val = "40.26" }
    @{ Addr = "D48"; Val = "127.32" }
    @{ Addr = "D50"; Val = "2.26" }
)

foreach ($item in $textCells) {
    $rng = $ws.Range($item.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Val
    $rng.Style = "Normal"
}

# Remaining cells (names, links, percentage strings) are not ambiguous
# with numbers, so they can be assigned directly.

$ws.Range("D2").Value = "64.402.74"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "3.158.45"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.150.95"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").Value = "  +3.86%  "
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").Value = "3.686.78"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.170.63"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "64.069.69"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  +9.22%  "
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +12.20%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E31").Value = "  +6.57%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("D35").Value = "0.0₃0854"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("E42").Value = "  +4.89%  "
$ws.Range("E43").Value = "  +6.10%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "2.929.17"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("E46").Value = "  +15.37%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("E50").Value = "  +3.38%  "
$ws.Range("E51").Value = "  -0.30%  "
